$d = $word.ActiveDocument

# Locate the paragraph that currently holds the split ":" + URL runs
# (styled blue/underlined, sz 32) for the lab2 video link.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd() -like "*HWMs9OfojtE*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate target paragraph"
}

# Replace that single paragraph with four paragraphs:
#  1) the original ":" + URL text merged into one run, with vi-VN lang
#     added to both the paragraph mark formatting and kept on the run
#  2) a new empty paragraph (blue/underline formatting, vi-VN lang)
#  3) a new paragraph with the "Link video lab3_KOT104:" text
#  4) a new paragraph with the new lab3 video URL (blue/underline, vi-VN)
$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:color w:val="4472C4" w:themeColor="accent1"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/><w:lang w:val="vi-VN"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:color w:val="4472C4" w:themeColor="accent1"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr><w:t>:https://youtu.be/HWMs9OfojtE?si=i7Mfh0tgK9kfzj6D</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p><w:pPr><w:rPr><w:color w:val="4472C4" w:themeColor="accent1"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/><w:lang w:val="vi-VN"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="vi-VN"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="vi-VN"/></w:rPr><w:t>Link video lab3_KOT104:</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="4472C4" w:themeColor="accent1"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/><w:lang w:val="vi-VN"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="4472C4" w:themeColor="accent1"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/><w:lang w:val="vi-VN"/></w:rPr><w:t>https://youtu.be/NkyiGxkalQY?si=SI1AYXSp_6mWqoVo</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$target.Range.InsertXML($xml)
